$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inst")

# --- Workbook-level bits -----------------------------------------------
$wb.Windows.Item(1).WindowState = $wb.Windows.Item(1).WindowState
$excel.Windows.Item(1).Left = 2145

# --- Insert a new first column (A) that holds a computed "Name" -------
$ws.Columns.Item(1).Insert()

# --- Insert a new "Action" column before the trailing Comment column --
# After the first insert, the old column L (Comment) is now M; insert
# before it so the new column becomes M and Comment shifts to N.
$ws.Columns.Item(13).Insert()

$wb.Save()
